# Rename variable codes so they zero-pad numerically (Rx1 -> Rx01, Ca1 -> Ca01, Ca1a -> Ca01a, ...)
# then re-sort the table by the "Variable #" column so the new names land in proper numeric order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 109

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()

    if ($val -eq "Rx1")  { $cell.Value = "Rx01" }
    elseif ($val -eq "Rx2")  { $cell.Value = "Rx02" }
    elseif ($val -eq "Rx3")  { $cell.Value = "Rx03" }
    elseif ($val -eq "Rx4")  { $cell.Value = "Rx04" }
    elseif ($val -eq "Rx5")  { $cell.Value = "Rx05" }
    elseif ($val -eq "Rx6")  { $cell.Value = "Rx06" }
    elseif ($val -eq "Rx7")  { $cell.Value = "Rx07" }
    elseif ($val -eq "Rx8")  { $cell.Value = "Rx08" }
    elseif ($val -eq "Rx9")  { $cell.Value = "Rx09" }
    elseif ($val -eq "Ca1")  { $cell.Value = "Ca01" }
    elseif ($val -eq "Ca1a") { $cell.Value = "Ca01a" }
    elseif ($val -eq "Ca2")  { $cell.Value = "Ca02" }
    elseif ($val -eq "Ca3")  { $cell.Value = "Ca03" }
    elseif ($val -eq "Ca4")  { $cell.Value = "Ca04" }
    elseif ($val -eq "Ca5")  { $cell.Value = "Ca05" }
    elseif ($val -eq "Ca6")  { $cell.Value = "Ca06" }
    elseif ($val -eq "Ca7")  { $cell.Value = "Ca07" }
    elseif ($val -eq "Ca8")  { $cell.Value = "Ca08" }
    elseif ($val -eq "Ca9")  { $cell.Value = "Ca09" }
}

# Re-sort the data range (excluding header) on column A, ascending, as the table's
# sort button would do after the renames shift the natural text order.
$sortRange = $ws.Range("A2:E$lastRow")
$keyRange  = $ws.Range("A2:A$lastRow")
$sortRange.Sort($keyRange)

# Leave the selection where Excel would land after performing the sort from the
# bottom of the visible window.
[void]$ws.Range("A4").Select()
